# Add tutorial strings ("showcase_drag_subject" / "showcase_tap_subject")
# to the localization table on the "string" sheet.
#
# The table is sorted alphabetically by the "keys" column, and the two new
# keys sort right after "setup" (row 114) and before "sort_by" (old row 115),
# so we insert two new rows at sheet row 115 and push everything below it
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert two blank rows right before the old row 115 ("sort_by"), shifting
# every row below it (up to the old last row, 147) down by two.
$ws.Rows("115:116").Insert()

# Grow the table/list-object definition (and its autofilter) to cover the
# two new rows, so the table keeps spanning the full data range (A1:E149).
$lo.Resize($ws.Range("A1:E149"))

# Fill in the new row for "showcase_drag_subject" (keys first, then the
# translated values in en/de/fr/lb order, matching the source commit).
$ws.Range("A115").Value = "showcase_drag_subject"
$ws.Range("A116").Value = "showcase_tap_subject"

$ws.Range("B115").Value = "Drag to change subject order"
$ws.Range("D115").Value = "Ziehe, um die Reihenfolge der Fächer zu ändern"
$ws.Range("C115").Value = "Glissez pour changer l'ordre des matières"
$ws.Range("E115").Value = "Zéi fir d'Reiefolleg vun de Fächer ze änneren"

$ws.Range("B116").Value = "Tap to make the subject above a subject group"
$ws.Range("C116").Value = "Appuyez pour faire du sujet au-dessus un groupe de matières"
$ws.Range("D116").Value = "Berühre um aus dem Fach drüber eine Fachgruppe zu machen"
$ws.Range("E116").Value = "//TODO"

# Match the row height/custom-height formatting used by every other data
# row in the sheet (the plain row Insert above leaves it unset).
$ws.Rows("115:116").RowHeight = 18.75

# Update the ExternalData_1 defined name so it keeps tracking the full
# imported range after the table grew by two rows.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*ExternalData_1*") {
        $n.RefersTo = "=string!`$A`$1:`$B`$149"
    }
}

# Restore the window scroll position / active cell selection used in the
# edited workbook.
$excel.ActiveWindow.ScrollRow = 127
$ws.Range("A131").Select()
